# Extend the "Eng actual exchange rate" sheet with two new daily columns:
#   DK -> 14.06.24
#   DL -> 17.06.24
# This mirrors the existing layout: columns are added by copying the
# formatting of the last populated column (DJ) and then filling in the
# new header labels, raw data and the "change vs previous day" formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting (styles/number formats/borders) of the last
#    column (DJ, rows 1-59) onto the two new columns so every new cell
#    gets the same look the sheet author used for every previous column.
$srcFmt = $ws.Range("DJ1:DJ59")
$dstFmt = $ws.Range("DK1:DL59")
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the "bestFit" width (16 characters rendered) that DA:DJ already use
# for this same data block.
$ws.Range("DK1:DL59").EntireColumn.ColumnWidth = 15.1666666666667

# 2) Header row (row 5) - new date labels.
$ws.Range("DK5").Value2 = "14.06.24"
$ws.Range("DL5").Value2 = "17.06.24"

# 3) Raw data rows (7-29) - values taken from the source workbook.
$rawData = @{
    7  = @(942500.65823890967, 939895.1588591066)
    8  = @(745918.73152969009, 743611.59776629356)
    9  = @(-392458.82517259999, -390767.90789720003)
    10 = @(670295.5266173, 669332.19195340003)
    11 = @(368164.106768, 368410.26691900002)
    12 = @(368164.106768, 368410.26691900002)
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @(0, 0)
    17 = @(0, 0)
    18 = @(0, 0)
    19 = @(0, 0)
    20 = @(0, 0)
    21 = @(0, 0)
    22 = @(468082.0300849902, 465047.31371009361)
    23 = @(1688419.3897685998, 1683506.7566254002)
    24 = @(889639.67366279999, 890811.5800357)
    25 = @(328355.24917009997, 327319.11003029998)
    26 = @(448837.34673009999, 447558.86425559997)
    27 = @(21587.120205599815, 17817.202303800208)
    28 = @(19564.9055974, 15441.157375000001)
    29 = @(2022.2146082000002, 2376.0449288000004)
}

foreach ($row in $rawData.Keys) {
    $vals = $rawData[$row]
    $ws.Range("DK$row").Value2 = $vals[0]
    $ws.Range("DL$row").Value2 = $vals[1]
}

# 4) "Change vs previous day" formula rows (33-55). Each formula row maps
#    to a raw-data row 26 above it, and is either "cur - prev" or, for the
#    repo-instrument block (rows 43-47 <-> data rows 17-21), the negated
#    difference "-(cur - prev)".
$negatedRows = 17..21

for ($dataRow = 7; $dataRow -le 29; $dataRow++) {
    $formulaRow = $dataRow + 26
    if ($negatedRows -contains $dataRow) {
        $ws.Range("DK$formulaRow").Formula = "=-(DK$dataRow-DJ$dataRow)"
        $ws.Range("DL$formulaRow").Formula = "=-(DL$dataRow-DK$dataRow)"
    } else {
        $ws.Range("DK$formulaRow").Formula = "=DK$dataRow-DJ$dataRow"
        $ws.Range("DL$formulaRow").Formula = "=DL$dataRow-DK$dataRow"
    }
}

# 5) Keep the sheet's current selection consistent with the new extent.
$ws.Range("DT26").Select()

$excel.CutCopyMode = 0
